# Fixed DIP and updated version number.
# Appends 18 new case rows (1153-1170) to Sheet1 of the Case_Data workbook,
# all for case 21TRC08418 / Hemmeter, mirroring the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value to a cell while forcing it to be stored as TEXT
# (matches the source data, where numeric-looking values such as case
# numbers, statute codes, dollar amounts and day counts are all text).
# NumberFormat "@" keeps Excel from re-interpreting the string as a number,
# and resetting the Style back to "Normal" afterwards drops the explicit
# style index again so the cell is left on the default (unstyled) format,
# just like every other cell already in this sheet.
function Set-TextCell {
    param($Row, $Col, $Value)
    $c = $ws.Cells.Item($Row, $Col)
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = "Normal"
}

$rows = @(
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="177" },
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="177" },
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="177" },
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="185" },
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";                 H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="177" },
    @{ A="21TRC08418"; B="Hemmeter"; C="Driving In Marked Lanes"; D="4511.33";     E="MM"; F="Dismissed";     G="";      H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="Turn And Stop Signals";   D="No Data";     E="MM"; F="Dismissed";     G="";      H=" ";     I=" ";   J=" ";  K=" "   },
    @{ A="21TRC08418"; B="Hemmeter"; C="OVI Alcohol / Drugs 1st"; D="4511.19A1A*"; E="M1"; F="Guilty";     G="Guilty";    H="$ 375"; I="$ 0"; J="180"; K="177" }
)

$startRow = 1153
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    Set-TextCell $r 1  $data.A
    Set-TextCell $r 2  $data.B
    Set-TextCell $r 3  $data.C
    Set-TextCell $r 4  $data.D
    Set-TextCell $r 5  $data.E
    Set-TextCell $r 6  $data.F
    if ($data.ContainsKey("G")) {
        Set-TextCell $r 7 $data.G
    }
    Set-TextCell $r 8  $data.H
    Set-TextCell $r 9  $data.I
    Set-TextCell $r 10 $data.J
    Set-TextCell $r 11 $data.K
}
